$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cloghan, Entekra COE", "Monaghan"),
    @("Derrylin O'Connells GAA, Derrylin", "Fermanagh"),
    @("Lemybrien", "Waterford"),
    @("Clontibret", "Monaghan"),
    @("Páirc Tailteann", "Meath"),
    @("LIT Gaelic Grounds, Limerick", "Limerick")
)

$row = 94
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

$ws.Range("A101").Select()
